$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: add the 7th test case (row 8) - "Not Null Columns" ---
$ws1.Range("A8").Value = 7
$ws1.Range("B8").Value = "Not Null Columns"
$ws1.Range("C8").Value = "CHECK_NOT_NULL:superstore_data"
$ws1.Range("D8").Value = '["Order ID","Customer ID","Product ID"]'
$ws1.Range("A8:D8").HorizontalAlignment = -4131

# Move Sheet1 selection to C7
[void]$ws1.Range("C7").Select()

# --- Sheet2: now holds only the 6th & 7th test cases ---
# Row 1: 6th test case - Schema Validation (query changed to the SQL variant, expected value updated)
$ws2.Range("A1").Value = 6
$ws2.Range("C1").Value = "SELECT column_name, data_type FROM information_schema.columns WHERE table_name = 'superstore_data';"
$ws2.Range("D1").Value = '{"order_id":"text","order_date":"text","region":"text","category":"text","profit":"double_precision"}'
$ws2.Range("C1").WrapText = $true
$ws2.Range("C1").HorizontalAlignment = -4131

# Row 2: 7th test case - Not Null Columns (expected value updated to new column names)
$ws2.Range("A2").Value = 7
$ws2.Range("D2").Value = '["Order ID","Customer ID","Product ID"]'

# Rows 3-7: clear old test-case data (styles/formatting remain untouched)
$ws2.Range("A3:D7").ClearContents()

# Move Sheet2 selection to A2:H2 with active cell A2
[void]$ws2.Range("A2:H2").Select()

# Sheet1 remains the active/visible tab
[void]$ws1.Activate()

Write-Output "done"
